$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 730.6461864673427
$ws.Range("C4").Value = 45.55891441616916
$ws.Range("C5").Value = 2860.383513718284
$ws.Range("D7").Value = 691.7720860150996
$ws.Range("D8").Value = 526.9099049891718
